$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-8 down to 5-9.
$ws.Rows.Item(4).Insert()

# Copy the style (incl. number format) of the date cell from the row above
# (row 5, which used to be row 4) onto the new D4 cell.
$ws.Range("D5").Copy()
$ws.Range("D4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new weekly record in row 4.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44589
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112001
$ws.Cells.Item(4, 7).Value = "Berenjena"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 110
$ws.Cells.Item(4, 11).Value = 5000
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = 5500
$ws.Cells.Item(4, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 92
$ws.Cells.Item(4, 17).Value = 60
$ws.Cells.Item(4, 18).Value = "Hortaliza"
